# Update the "dSF" column (F) values on the specific rows that were
# re-pulled / recalculated, per the commit "repull data, push all data,
# mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    5  = -2
    13 = 3
    17 = -3
    20 = 5
    23 = -4
    26 = -7
    29 = 3
    30 = 2
    36 = -4
    39 = -2
    42 = 1
    45 = 2
    46 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
